# Update cryptos list row-by-row per the latest scrape.
# All data cells in this sheet are stored as text (inline strings),
# so we force NumberFormat = "@" before writing each value to keep
# Excel from auto-converting numeric-looking text to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.393.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.40%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.67'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.13%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '264.63'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.01%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5209'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.88%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3267'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.82%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06803'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.19%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.84'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.80%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7746'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.06%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07774'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.17%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.811.66'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.81%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.92'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.011'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.03%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.16%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.91'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.29%  '

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.06%  '

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007977'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.45%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.406.61'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.47%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.071.74'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.629'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.52%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.572'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.986'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.38%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.29'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.18%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.178'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -9.22%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.667'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.83%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.43%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.07'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.71%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.168'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.63%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.130'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.79%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08741'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.30%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04820'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.73%  '

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.134'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.11%  '

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7204'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.89%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.856'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.63%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.093'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.17%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01778'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.209'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.38%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4858'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.28%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9109'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.49%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '111.01'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.29%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.054'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.43%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.705'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.52%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05934'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.06%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4159'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.59%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.077'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.69%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1236'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.65%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.93'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.08%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8870'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.18%  '

